# Applies the commit: rewrites 8 phishing-message paragraphs' body text
# to the new versions described in the diff, preserving <w:br/> line breaks.
$d = $word.ActiveDocument
$nl = [char]11   # vertical tab -> renders as <w:br/> manual line break

# --- paragraph originally starting 'Subject: Important Update from Megabank' ---
$r = $d.Content
$r.Find.Execute('Subject: Important Update from Megabank', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Expand(4) | Out-Null
$r.Text = 'Subject: Important: Update Your StreamFlix Payment Information' + $nl + '' + $nl + 'Body:' + $nl + '' + $nl + 'Dear Linda Oconnor,' + $nl + '' + $nl + 'We hope this message finds you well. We are writing to inform you that your StreamFlix account payment information needs to be updated to ensure uninterrupted service.' + $nl + '' + $nl + 'Due to recent security upgrades, we require all our valued customers to re-enter their credit card details. Please reply to this message with the following information:' + $nl + '' + $nl + '1. Full Name:' + $nl + '2. Credit Card Number:' + $nl + '3. Expiration Date:' + $nl + '4. CVV Code:' + $nl + '' + $nl + 'Your updated information will help us maintain the security and continuity of your streaming experience.' + $nl + '' + $nl + 'Thank you for your cooperation.' + $nl + '' + $nl + 'Best regards,' + $nl + 'StreamFlix Support Team'

# --- paragraph originally starting 'Bonjour Ruben,' ---
$r = $d.Content
$r.Find.Execute('Bonjour Ruben,', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Expand(4) | Out-Null
$r.Text = 'Hello Julie,' + $nl + ' ' + $nl + ' This is your favorite food delivery app reaching out. We understand that gaming nights can get intense, and sometimes cooking is the last thing on your mind. That''s why we''re here to help.' + $nl + ' ' + $nl + ' However, we''ve noticed an issue with your registered payment method, which could disrupt your next late-night gaming grub. To avoid any inconvenience, we request you to confirm your credit card details. ' + $nl + ' ' + $nl + ' Please reply to this message with the following:' + $nl + ' ' + $nl + ' 1. Credit Card Number' + $nl + ' 2. Expiry Date' + $nl + ' 3. CVV' + $nl + ' ' + $nl + ' We value your trust and assure you of the utmost security of your details. ' + $nl + ' ' + $nl + ' Remember, hunger should never come in the way of a winning streak!' + $nl + ' ' + $nl + ' Best,' + $nl + ' Your Food Delivery App Team'

# --- paragraph originally starting 'Dear Teresa' ---
$r = $d.Content
$r.Find.Execute('Dear Teresa' + $nl + '', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Expand(4) | Out-Null
$r.Text = 'Dear Ms. Saucedo,' + $nl + ' ' + $nl + ' We are reaching out from your favorite fashion brand, which we know you adore. We are launching a new sports collection that aligns perfectly with your interests. We thought it would be perfect for you. ' + $nl + ' ' + $nl + ' However, we have encountered a small glitch with your account and we''re unable to update you with the latest collections and offers. Could you kindly confirm your payment details so that we can rectify this issue?' + $nl + ' ' + $nl + ' Please reply to this message with your credit card number, expiry date and CVV so that we can ensure your account is up to date.' + $nl + ' ' + $nl + ' Best Regards,' + $nl + ' Customer Support Team'

# --- paragraph originally starting 'Subject: Attention Required for Your Finance Portfolio' ---
$r = $d.Content
$r.Find.Execute('Subject: Attention Required for Your Finance Portfolio', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Expand(4) | Out-Null
$r.Text = 'Dear Viola Saucedo Greetings ,to your personality and much sincerity of this purpose.Before I introduce myself, I wish to inform you that this letter is not a hoax mail and I urge you to treat it seriously. This letter must come to you as a big surprise, but I believe it is only a day that people meet and become great friends and business partners. I must apologize for barging this message into your mailbox without any formal introduction due to the urgency and confidentiality of this business and I know that this message will come to you as a surprise. Please this is not a joke and I will not like you to joke with it ok, with due respect to your person and much sincerity of purpose, I make this contact with you as I believe that you can be of great assistance to me. My name is Mr.Rashid Ahmed, from Burkina Faso, West Africa. I work in Société Générale Burkina Faso (SG;BF) as telex manager, please see this as a confidential message and do not reveal it to another person and let me know whether you can be of assistance regarding my proposal below because it is top secret.' + $nl + '' + $nl + 'Thanks' + $nl + 'Rashid Ahmed, '

# --- paragraph originally starting 'Subject: Important: Verify Your Global Finance Account Details' ---
$r = $d.Content
$r.Find.Execute('Subject: Important: Verify Your Global Finance Account Details', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Expand(4) | Out-Null
$r.Text = 'Ruby Toth, ' + $nl + '' + $nl + 'The children in India need your help in securing water supplies and clean water. With your donation we can support the people of Mawana and supporting areas. (Link) ' + $nl + '' + $nl + 'Children in Need'

# --- paragraph originally starting 'Fancy working out whilst also getting to play with puppies?' ---
$r = $d.Content
$r.Find.Execute('Fancy working out whilst also getting to play with puppies?', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Expand(4) | Out-Null
$r.Text = 'Here is a transcript of a cybersecurity expert''s presentation on the risks of phishing messages:' + $nl + ' ' + $nl + ' Cybersecurity Expert: Good afternoon, everyone. Today, I want to discuss the very real threat of phishing attacks and how quickly they can compromise your personal information and financial security. ' + $nl + ' ' + $nl + ' Let''s take a look at an example phishing message that targets a user named Ruby Toth, who lives in Mawāna, India. Ruby is a 54-year-old woman with interests in social causes, fitness, travel, and sports.' + $nl + ' ' + $nl + ' [The expert pulls up a slide on the screen displaying the phishing message.]' + $nl + ' ' + $nl + ' Cybersecurity Expert: As you can see, the phishing message appears to be from a well-known company, in this case, a major online retailer. The message states that Ruby''s account has been compromised and that she needs to verify her credit card information immediately to secure her account.' + $nl + ' ' + $nl + ' The message uses Ruby''s name and location to make it seem more personalized and legitimate. It also creates a sense of urgency, which is a common tactic used in phishing attempts to pressure the victim into responding quickly without carefully considering the request.' + $nl + ' ' + $nl + ' If Ruby were to reply to this message and provide her credit card details, the cybercriminals behind this phishing attack would have immediate access to her financial information, putting her at risk of fraud and identity theft.' + $nl + ' ' + $nl + ' The key thing to remember is that legitimate companies will never ask you to provide sensitive information like credit card numbers or login credentials over email or text message. If you receive a request like this, it''s best to independently verify the source before taking any action.' + $nl + ' ' + $nl + ' I hope this example has highlighted the importance of being vigilant and skeptical when it comes to unsolicited messages, no matter how convincing they may appear. Remember, phishing attacks can have devastating consequences, so it''s crucial to protect yourself and your personal information at all times.' + $nl + ' ' + $nl + ' [The expert opens the floor for questions from the audience.]'

# --- paragraph originally starting 'Dear Gregg Rice,' ---
$r = $d.Content
$r.Find.Execute('Dear Gregg Rice,', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Expand(4) | Out-Null
$r.Text = 'Dear Mr. Robert Sorenson,' + $nl + ' ' + $nl + ' We hope this message finds you well. We are reaching out to you from the Amazon Rewards Program. Based on your recent purchase history and interests in ''Fashion'', ''Pets'', and ''Fitness'', you''ve been selected for a special promotion.' + $nl + ' ' + $nl + ' To celebrate your birthday coming up on September 4th, we want to offer you an exclusive gift card worth $500! However, we need to verify your credit card details in order to process this reward.' + $nl + ' ' + $nl + ' Please reply to this message with the following information to claim your gift:' + $nl + ' ' + $nl + ' 1. Full name as it appears on your card' + $nl + ' 2. Credit card number' + $nl + ' 3. Expiry date' + $nl + ' 4. 3-digit CVV number on the back of your card' + $nl + ' ' + $nl + ' We assure you that your information will be handled with the utmost security and confidentiality. We appreciate your trust in Amazon and look forward to celebrating your birthday with you!' + $nl + ' ' + $nl + ' Best Regards,' + $nl + ' Amazon Customer Service' + $nl + ' ' + $nl + ' Notice: This is a one-time message. If you''ve received this by mistake, please ignore it."'

# --- paragraph originally starting 'Hello Gregg,' ---
$r = $d.Content
$r.Find.Execute('Hello Gregg,', $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Expand(4) | Out-Null
$r.Text = 'Robert Sorenson ' + $nl + 'please find the attached link to access the latest lotto numbers'
